$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => (new E value, new F value or $null if unchanged)
$updates = @(
    @(2, 9, $null),
    @(3, 9, $null),
    @(4, 9, $null),
    @(5, 7, $null),
    @(6, 9, $null),
    @(7, 7, $null),
    @(8, 9, $null),
    @(9, 7, $null),
    @(10, 2, $null),
    @(11, 9, $null),
    @(12, 7, $null),
    @(13, 9, $null),
    @(14, 9, $null),
    @(15, 9, $null),
    @(16, 1, $null),
    @(17, 7, $null),
    @(18, 10, 20251108),
    @(19, 10, 20251108),
    @(20, 10, 20251108),
    @(21, 10, 20251108),
    @(22, 7, $null),
    @(23, 7, $null),
    @(24, 7, $null),
    @(25, 7, $null),
    @(26, 7, $null),
    @(27, 3, $null),
    @(28, 10, 20251108),
    @(29, 10, 20251108),
    @(30, 10, 20251108),
    @(31, 10, 20251108),
    @(32, 10, 20251108),
    @(33, 10, 20251108),
    @(34, 10, 20251108),
    @(35, 10, 20251108),
    @(37, 10, 20251108),
    @(38, 10, 20251108),
    @(39, 10, 20251108),
    @(40, 2, $null),
    @(41, 2, $null),
    @(42, 10, 20251108),
    @(43, 7, $null),
    @(44, 2, $null),
    @(45, 7, $null),
    @(46, 2, $null),
    @(47, 10, 20251108),
    @(48, 2, $null),
    @(49, 3, $null),
    @(50, 5, $null),
    @(51, 5, $null),
    @(52, 5, $null),
    @(53, 5, $null),
    @(54, 5, $null),
    @(55, 5, $null),
    @(56, 5, $null),
    @(57, 5, $null),
    @(58, 9, $null),
    @(59, 9, $null),
    @(60, 9, $null),
    @(61, 3, $null),
    @(62, 9, $null),
    @(63, 9, $null),
    @(64, 9, $null),
    @(65, 10, 20251108),
    @(66, 10, 20251108),
    @(67, 10, 20251108),
    @(68, 10, 20251108),
    @(69, 10, 20251108),
    @(70, 1, $null),
    @(71, 1, $null),
    @(72, 1, $null),
    @(73, 1, $null),
    @(74, 1, $null),
    @(75, 1, $null),
    @(76, 1, $null),
    @(77, 4, $null),
    @(78, 4, $null),
    @(79, 4, $null),
    @(80, 4, $null),
    @(81, 4, $null),
    @(82, 4, $null),
    @(83, 4, $null),
    @(84, 4, $null),
    @(85, 4, $null),
    @(86, 4, $null),
    @(87, 2, $null),
    @(88, 2, $null),
    @(89, 2, $null),
    @(90, 2, $null),
    @(91, 7, $null),
    @(92, 2, $null),
    @(93, 4, $null),
    @(94, 5, $null),
    @(95, 3, $null),
    @(96, 1, $null),
    @(97, 1, $null),
    @(98, 1, $null),
    @(99, 1, $null),
)

foreach ($u in $updates) {
    $row = $u[0]
    $eVal = $u[1]
    $fVal = $u[2]
    $ws.Cells.Item($row, 5).Value = $eVal
    if ($fVal -ne $null) {
        $ws.Cells.Item($row, 6).Value = $fVal
    }
}
